$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.727.69"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.253.41"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.78"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.53"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.444"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0970"
$ws.Range("E10").Value = "  -6.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.55"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.51"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "2.587.12"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.60"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.11"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").Value = "2.261.93"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "43.661.86"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "0.0₃0975"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.59"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.26"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.66"
$ws.Range("E25").Value = "  +31.50%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.22"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.127"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.92"
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0681"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.93"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.36"
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("E39").Value = "  -2.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0254"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.59"
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.63"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.20"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.16"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.18"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0942"
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("D48").Value = "1.451.81"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000208"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.84"
$ws.Range("E51").Value = "  -6.62%  "
